$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.555.25'
$ws.Range("E2").Value = '  -2.67%  '
$ws.Range("D3").Value = '1.860.89'
$ws.Range("E3").Value = '  -2.56%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '286.86'
$ws.Range("E5").Value = '  -6.86%  '
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5221'
$ws.Range("E7").Value = '  -2.72%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3685'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07082'
$ws.Range("E9").Value = '  -3.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.03'
$ws.Range("E10").Value = '  -4.66%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8752'
$ws.Range("E11").Value = '  -3.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08038'
$ws.Range("E12").Value = '  -2.03%  '
$ws.Range("D13").Value = '1.891.22'
$ws.Range("E13").Value = '  +65.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.226'
$ws.Range("E14").Value = '  -2.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '90.62'
$ws.Range("E15").Value = '  -5.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.003'
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.53'
$ws.Range("E17").Value = '  -2.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008398'
$ws.Range("E18").Value = '  -2.95%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.004'
$ws.Range("E19").Value = '  +0.14%  '
$ws.Range("D20").Value = '26.590.12'
$ws.Range("E20").Value = '  -2.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.911'
$ws.Range("E21").Value = '  -2.74%  '
$ws.Range("E22").Value = '  -2.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.284'
$ws.Range("E23").Value = '  -3.63%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '144.57'
$ws.Range("E24").Value = '  -3.51%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.231'
$ws.Range("E25").Value = '  -2.73%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.732'
$ws.Range("E26").Value = '  -0.84%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.79'
$ws.Range("E27").Value = '  -2.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '112.70'
$ws.Range("E28").Value = '  -3.68%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.630'
$ws.Range("E29").Value = '  -4.25%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.541'
$ws.Range("E30").Value = '  -5.46%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08999'
$ws.Range("E31").Value = '  -3.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7823'
$ws.Range("E32").Value = '  -6.45%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04924'
$ws.Range("E33").Value = '  -2.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.150'
$ws.Range("E34").Value = '  -6.22%  '
$ws.Range("E35").Value = '  -3.20%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.5859'
$ws.Range("E36").Value = '  +1.53%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01917'
$ws.Range("E39").Value = '  -4.64%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.036'
$ws.Range("E40").Value = '  -3.97%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.408'
$ws.Range("E41").Value = '  -2.46%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.559'
$ws.Range("E44").Value = '  -8.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1468'
$ws.Range("E45").Value = '  -3.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.005'
$ws.Range("E46").Value = '  +0.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.839'
$ws.Range("E47").Value = '  -3.03%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.596'
$ws.Range("E48").Value = '  -2.70%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.74'
$ws.Range("E49").Value = '  -4.66%  '
$ws.Range("E50").Value = '  -1.97%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '61.48'
$ws.Range("E51").Value = '  -3.03%  '

# Row 37 and 38 swap (MXToken <-> RenderToken)
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.590'
$ws.Range("E37").Value = '  -4.55%  '

$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.156'
$ws.Range("E38").Value = '  -5.96%  '

# Row 42 and 43 swap (Quant <-> Decentraland)
$ws.Range("B42").Value = 'Decentraland'
$ws.Range("C42").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5111'
$ws.Range("E42").Value = '  +3.63%  '

$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '114.00'
$ws.Range("E43").Value = '  -3.08%  '
